$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1585.9
$ws.Range("I32").Value = 1166
$ws.Range("J32").Value = 1865.8334
$ws.Range("K32").Value = 1166
$ws.Range("L32").Value = 1865.8334
$ws.Range("M32").Value = -840
$ws.Range("N32").Value = -2517.8334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 205.33333
$ws.Range("I42").Value = 58
$ws.Range("K42").Value = 174
$ws.Range("M42").Value = 56

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 715805.9
$ws.Range("I80").Value = 1615.75
$ws.Range("J80").Value = 1001481.9
$ws.Range("K80").Value = 4847.25
$ws.Range("L80").Value = 3004445.7
$ws.Range("M80").Value = -3849.25
$ws.Range("N80").Value = -3006441.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 715805.9
$ws.Range("I83").Value = 1615.75
$ws.Range("J83").Value = 1001481.9
$ws.Range("K83").Value = 14541.75
$ws.Range("L83").Value = 9013337.1
$ws.Range("M83").Value = -9549.75
$ws.Range("N83").Value = -9023321.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 966.9
$ws.Range("I103").Value = 884.5
$ws.Range("J103").Value = 987.5
$ws.Range("K103").Value = 2653.5
$ws.Range("L103").Value = 2962.5
$ws.Range("M103").Value = -2067.5
$ws.Range("N103").Value = -4134.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2402.0293
$ws.Range("I112").Value = 988.5
$ws.Range("J112").Value = 2490.375
$ws.Range("K112").Value = 2965.5
$ws.Range("L112").Value = 7471.125
$ws.Range("M112").Value = -1857.5
$ws.Range("N112").Value = -9687.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1234.9487
$ws.Range("I132").Value = 1248.5714
$ws.Range("J132").Value = 1115.75
$ws.Range("K132").Value = 3745.7142
$ws.Range("L132").Value = 3347.25
$ws.Range("M132").Value = -1215.7142
$ws.Range("N132").Value = -8407.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 12476
$ws.Range("I137").Value = 5465.3125
$ws.Range("J137").Value = 17817.477
$ws.Range("K137").Value = 16395.9375
$ws.Range("L137").Value = 53452.431
$ws.Range("M137").Value = -13845.9375
$ws.Range("N137").Value = -58552.431

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2235.25
$ws.Range("I138").Value = 1805.6111
$ws.Range("J138").Value = 2664.889
$ws.Range("K138").Value = 5416.8333
$ws.Range("L138").Value = 7994.667
$ws.Range("M138").Value = -276.8333000000002
$ws.Range("N138").Value = -18274.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 26054.5
$ws.Range("I74").Value = 29028.615
$ws.Range("K74").Value = 29028.615
$ws.Range("M74").Value = -28154.615

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 26054.5
$ws.Range("I77").Value = 29028.615
$ws.Range("K77").Value = 145143.075
$ws.Range("M77").Value = -140775.075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 283.42856
$ws.Range("I22").Value = 333.75
$ws.Range("K22").Value = 333.75
$ws.Range("M22").Value = -160.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6871.85
$ws.Range("I134").Value = 2946.756
$ws.Range("J134").Value = 15341.789
$ws.Range("K134").Value = 8840.268
$ws.Range("L134").Value = 46025.367
$ws.Range("M134").Value = -6305.268
$ws.Range("N134").Value = -51095.367

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4250.96
$ws.Range("I31").Value = 4375.778
$ws.Range("K31").Value = 4375.778
$ws.Range("M31").Value = -4080.778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4250.96
$ws.Range("I34").Value = 4375.778
$ws.Range("K34").Value = 4375.778
$ws.Range("M34").Value = -4173.778

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 35991
$ws.Range("J44").Value = 3987
$ws.Range("L44").Value = 11961
$ws.Range("N44").Value = -12757

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 1999.5
$ws.Range("J49").Value = 1999.5
$ws.Range("L49").Value = 5998.5
$ws.Range("N49").Value = -6310.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3332.3333
$ws.Range("I64").Value = 2999
$ws.Range("J64").Value = 3999
$ws.Range("K64").Value = 8997
$ws.Range("L64").Value = 11997
$ws.Range("M64").Value = -8727
$ws.Range("N64").Value = -12537

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 3332.3333
$ws.Range("I67").Value = 2999
$ws.Range("J67").Value = 3999
$ws.Range("K67").Value = 8997
$ws.Range("L67").Value = 11997
$ws.Range("M67").Value = -8061
$ws.Range("N67").Value = -13869

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1038.5454
$ws.Range("J113").Value = 763.1429000000001
$ws.Range("L113").Value = 2289.4287
$ws.Range("N113").Value = -6629.4287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 6667321.5
$ws.Range("I122").Value = 383.5
$ws.Range("J122").Value = 11111946
$ws.Range("K122").Value = 3451.5
$ws.Range("L122").Value = 100007514
$ws.Range("M122").Value = -1001.5
$ws.Range("N122").Value = -100012414

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2670.375
$ws.Range("J131").Value = 2942.0588
$ws.Range("L131").Value = 8826.1764
$ws.Range("N131").Value = -18906.1764

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 741.8570999999999
$ws.Range("I102").Value = 771.45
$ws.Range("K102").Value = 771.45
$ws.Range("M102").Value = 850.55

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3199.6924
$ws.Range("I122").Value = 3054.3635
$ws.Range("K122").Value = 9163.0905
$ws.Range("M122").Value = -6713.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8579.6875
$ws.Range("I7").Value = 6295.3335
$ws.Range("J7").Value = 9106.846
$ws.Range("K7").Value = 6295.3335
$ws.Range("L7").Value = 9106.846
$ws.Range("M7").Value = -6183.3335
$ws.Range("N7").Value = -9330.846

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2078.5625
$ws.Range("I22").Value = 566.3333
$ws.Range("K22").Value = 566.3333
$ws.Range("M22").Value = -271.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2078.5625
$ws.Range("I27").Value = 566.3333
$ws.Range("K27").Value = 566.3333
$ws.Range("M27").Value = -459.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1525.5385
$ws.Range("I46").Value = 990.087
$ws.Range("J46").Value = 2295.25
$ws.Range("K46").Value = 990.087
$ws.Range("L46").Value = 2295.25
$ws.Range("M46").Value = -802.087
$ws.Range("N46").Value = -2671.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 412.0625
$ws.Range("I55").Value = 291.18182
$ws.Range("K55").Value = 291.18182
$ws.Range("M55").Value = -118.18182

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2567.182
$ws.Range("I61").Value = 2574.9
$ws.Range("K61").Value = 2574.9
$ws.Range("M61").Value = -2372.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2250.111
$ws.Range("I82").Value = 1886.8
$ws.Range("K82").Value = 1886.8
$ws.Range("M82").Value = -1525.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2250.111
$ws.Range("I85").Value = 1886.8
$ws.Range("K85").Value = 1886.8
$ws.Range("M85").Value = -638.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2567.182
$ws.Range("I113").Value = 2574.9
$ws.Range("K113").Value = 2574.9
$ws.Range("M113").Value = -404.9000000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 8579.6875
$ws.Range("I126").Value = 6295.3335
$ws.Range("J126").Value = 9106.846
$ws.Range("K126").Value = 18886.0005
$ws.Range("L126").Value = 27320.538
$ws.Range("M126").Value = -16416.0005
$ws.Range("N126").Value = -32260.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 4973
$ws.Range("J74").Value = 2561.5
$ws.Range("L74").Value = 2561.5
$ws.Range("N74").Value = -4433.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 4973
$ws.Range("J77").Value = 2561.5
$ws.Range("L77").Value = 7684.5
$ws.Range("N77").Value = -17044.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 139407.73
$ws.Range("I132").Value = 253925.14
$ws.Range("K132").Value = 761775.42
$ws.Range("M132").Value = -759245.42
